$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "OLIVEIRA FRILLS"
$ws.Range("B3").Value = "https://res.cloudinary.com/maaji/image/upload/v1542895401/Spring2019/2110SBC08_2110SCC08_BLUE_1.jpg"
$ws.Range("C3").Value = "https://res.cloudinary.com/maaji/image/upload/v1542895404/Spring2019/2110SBC08_2110SCC08_BLUE_2.jpg"

$ws.Range("B7").Select()
